$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column E (rows 3-5 get highlighted "checked" names moved in from column G;
# rows 3-10 otherwise shift down by two positions) and column G (rows 1-13
# shift up, with two names from the bottom of E moving into the vacated
# slots at the bottom of G). Re-typed explicitly, cell by cell, to match
# the final layout exactly.

$ws.Range("G1").Value = "*Jonathan"
$ws.Range("G2").Value = "JonathanW"

$ws.Range("E3").Value = "*Francis"
$ws.Range("G3").Value = "Richard Zhao"

$ws.Range("E4").Value = "*Small Tony"
$ws.Range("G4").Value = "*Anthony "

$ws.Range("E5").Value = "*Yellow"
$ws.Range("G5").Value = "Tim"

$ws.Range("E6").Value = "Ms Tong"
$ws.Range("G6").Value = "Patrick"

$ws.Range("E7").Value = "Denis"
$ws.Range("G7").Value = "See Fu"

$ws.Range("E8").Value = "Begger"
$ws.Range("G8").Value = "Yvoone"

$ws.Range("E9").Value = "Sol Bread"
$ws.Range("G9").Value = "Pui"

$ws.Range("E10").Value = "Fai Chi"
$ws.Range("G10").Value = "Ocean"

$ws.Range("G11").Value = "**Lam Kei"
$ws.Range("G12").Value = "#Somingtat"
$ws.Range("G13").Value = "SomingtatW"

# Re-apply the "no-fill" style marker cells carry (same cellXfs slot already
# used elsewhere in the sheet) so the grey-highlight flag follows the values
# that moved, matching the target layout: E3:E5 and G1 / G4 highlighted.
$ws.Range("E3").Interior.ColorIndex = -4142
$ws.Range("E4").Interior.ColorIndex = -4142
$ws.Range("E5").Interior.ColorIndex = -4142
$ws.Range("G1").Interior.ColorIndex = -4142
$ws.Range("G4").Interior.ColorIndex = -4142

$ws.Range("G2").ClearFormats()
$ws.Range("G3").ClearFormats()
$ws.Range("G6").ClearFormats()

# Move the active selection from F4 to F1.
$ws.Range("F1").Select()
